# Regenerate merged AHB files:
#  - rename the "_old" / "_new" header suffixes to the concrete
#    form-version tags "_FV2404" / "_FV2410"
#  - turn the used range into a real Excel Table ("Table1")
#  - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row labels ------------------------------------------
# Columns A:J were the "<name>_old" block -> "<name>_FV2404"
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"

# Column K ("diff") is unchanged.

# Columns L:U were the "<name>_new" block -> "<name>_FV2410"
$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- 2) Turn A1:U66 into a Table with an autofilter -----------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), [System.Type]::Missing, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3) Freeze the header row ----------------------------------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
